# Expand the "tiers" sheet from 35 to 71 Philippine provinces (alphabetically
# sorted), each with a resilience and a risk category (Mid/Low/High).
# Policy assessment cards and more cleaning.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$provinces = @("Abra","Agusan Del Norte","Agusan Del Sur","Aklan","Albay","Antique","Apayao","Aurora","Basilan","Bataan","Batangas","Biliran","Bohol","Bukidnon","Bulacan","Cagayan","Camarines Norte","Camarines Sur","Camiguin","Capiz","Catanduanes","Cavite","Cebu","Compostela Valley","Davao Del Norte","Davao Del Sur","Davao Oriental","Eastern Samar","Guimaras","Ifugao","Ilocos Norte","Ilocos Sur","Iloilo","Isabela","Kalinga","La Union","Laguna","Lanao Del Norte","Lanao Del Sur","Leyte","Maguindanao","Marinduque","Masbate","Misamis Oriental","Negros Occidental","North Cotabato","Northern Samar","Nueva Ecija","Nueva Vizcaya","Occidental Mindoro","Oriental Mindoro","Palawan","Pampanga","Pangasinan","Quezon","Rizal","Romblon","Samar","Sarangani","Sorsogon","South Cotabato","Southern Leyte","Sultan Kudarat","Sulu","Surigao Del Norte","Tarlac","Tawi-Tawi","Zambales","Zamboanga Del Norte","Zamboanga Del Sur","Zamboanga Sibugay")
$resilience = @("Mid","Mid","Low","Mid","Mid","Low","Mid","High","Low","High","High","High","Mid","Low","High","High","Mid","Mid","Low","Mid","High","High","High","Low","Mid","High","Low","Low","Mid","Mid","High","High","High","Mid","High","High","High","Mid","Low","Mid","Low","Mid","Low","High","Mid","Low","Low","Mid","High","Mid","Mid","Mid","High","High","High","High","Low","Low","Low","Low","Mid","Low","Low","Low","Low","High","Low","High","Low","Mid","Low")
$risk = @("Mid","High","High","High","Low","High","Low","Mid","Mid","Low","Low","High","Mid","Low","Low","Mid","High","Mid","Mid","Mid","Mid","Mid","Mid","Low","Low","Low","Mid","High","High","Low","Low","Low","High","Mid","Low","Low","Low","Mid","Mid","High","High","High","High","High","Mid","Mid","High","Mid","Mid","High","High","High","Low","Mid","Low","Low","High","Mid","Low","High","Low","Mid","High","High","High","Mid","High","Low","High","Low","Mid")

# Existing rows are 2..36 (35 provinces); new rows 37..72 need the same
# formatting as column A already has for rows 2..36 (bold/bordered/centered
# style), so copy that format down before filling in values.
$ws.Range("A2").Copy()
$ws.Range("A37:A72").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

for ($i = 0; $i -lt $provinces.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $provinces[$i]
    $ws.Cells.Item($row, 2).Value = $resilience[$i]
    $ws.Cells.Item($row, 3).Value = $risk[$i]
}
